# Expand solutions for EC2 auth & deployment pipeline
# Resize/reposition several process boxes + connectors, and expand several
# box captions with the concrete tool names used for each pipeline stage.

function EMU($v) {
    # Shape.Left/Top/Width/Height are stored as single-precision points;
    # nudging the EMU value by +0.5 before the /12700 division keeps the
    # round-trip exact for the target EMU integer.
    return ($v + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Shape id=54 "Build application container or OS image" - taller box
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(1)
$sh.Left = EMU(3104514)
$sh.Top = EMU(1034410)
$sh.Width = EMU(1137330)
$sh.Height = EMU(794370)

# ---------------------------------------------------------------------
# 2) Shape id=55 "Continuous Event Monitoring & Threat Detection (...)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(2)
$sh.Left = EMU(7361500)
$sh.Top = EMU(4205921)
$sh.Width = EMU(1551006)
$sh.Height = EMU(759617)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Continuous Event Monitoring & Threat "
$tr.InsertAfter("Detection (")
$tr.InsertAfter("CloudWatch")
$tr.InsertAfter(", WAF, ")
$tr.InsertAfter("GuardDuty")
$tr.InsertAfter(")")

# ---------------------------------------------------------------------
# 3) Shape id=57 "Application code or OS change" - taller box
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(4)
$sh.Left = EMU(297060)
$sh.Top = EMU(1028060)
$sh.Width = EMU(1137330)
$sh.Height = EMU(800720)

# ---------------------------------------------------------------------
# 4) Shape id=58 "Infrastructure code change" - moved & taller
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(5)
$sh.Left = EMU(319019)
$sh.Top = EMU(3169350)
$sh.Width = EMU(1137330)
$sh.Height = EMU(821513)

# ---------------------------------------------------------------------
# 5) Connector id=59 (57 -> 54)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(6)
$sh.Left = EMU(1434390)
$sh.Top = EMU(1428420)
$sh.Width = EMU(1670124)
$sh.Height = EMU(3175)

# ---------------------------------------------------------------------
# 6) Connector id=60 (58 -> 56), flipV
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(7)
$sh.Left = EMU(1456349)
$sh.Top = EMU(2870550)
$sh.Width = EMU(4403298)
$sh.Height = EMU(709557)

# ---------------------------------------------------------------------
# 7) Connector id=61 (54 -> 56)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(8)
$sh.Left = EMU(4241844)
$sh.Top = EMU(1431595)
$sh.Width = EMU(1617803)
$sh.Height = EMU(841355)

# ---------------------------------------------------------------------
# 8) Connector id=62 (56 -> 55)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(9)
$sh.Left = EMU(6458673)
$sh.Top = EMU(2571750)
$sh.Width = EMU(1678330)
$sh.Height = EMU(1634171)

# ---------------------------------------------------------------------
# 9) Shape id=11 "Security & Compliance Scan of IAC Templates (Regula/Cloudsploit)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(10)
$sh.Left = EMU(2956135)
$sh.Top = EMU(3169350)
$sh.Width = EMU(1372798)
$sh.Height = EMU(821513)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Security & Compliance Scan of IAC "
$tr.InsertAfter("Templates")
$tr.InsertAfter([char]11)
$tr.InsertAfter("(Regula/")
$tr.InsertAfter("Cloudsploit")
$tr.InsertAfter(") ")

# ---------------------------------------------------------------------
# 10) Shape id=29 "Container / Image Scan (Clair)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(11)
$sh.Left = EMU(4508241)
$sh.Top = EMU(1028040)
$sh.Width = EMU(1137330)
$sh.Height = EMU(800740)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Container / Image "
$tr.InsertAfter("Scan (Clair)")

# ---------------------------------------------------------------------
# 11) Shape id=30 "Static Code Analysis & Security Scan" + "(Sonar/Snyk)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(12)
$sh.Left = EMU(1634188)
$sh.Top = EMU(1028040)
$sh.Width = EMU(1270528)
$sh.Height = EMU(800740)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Static Code Analysis & Security "
$tr.InsertAfter("Scan")
$tr.InsertAfter([char]13 + "(Sonar/")
$tr.InsertAfter("Snyk")
$tr.InsertAfter(")")

# ---------------------------------------------------------------------
# 12) Shape id=69 "Post Deployment Configuration Scan" + "(AWS Config, Inspector, SecurityHub)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(14)
$sh.Left = EMU(7361500)
$sh.Top = EMU(3169350)
$sh.Width = EMU(1551006)
$sh.Height = EMU(667686)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Post Deployment Configuration "
$tr.InsertAfter("Scan")
$tr.InsertAfter([char]13 + "(AWS ")
$tr.InsertAfter("Config")
$tr.InsertAfter(", Inspector, ")
$tr.InsertAfter("SecurityHub")
$tr.InsertAfter(")")
